$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(6, 8).Value = 105.75
$ws.Cells.Item(6, 9).Value = 106.181816
$ws.Cells.Item(6, 10).Value = 101
$ws.Cells.Item(6, 11).Value = 318.545448
$ws.Cells.Item(6, 12).Value = 303
$ws.Cells.Item(6, 13).Value = -206.545448
$ws.Cells.Item(6, 14).Value = -527

$ws.Cells.Item(12, 8).Value = 244
$ws.Cells.Item(12, 10).Value = 82
$ws.Cells.Item(12, 12).Value = 82
$ws.Cells.Item(12, 14).Value = -422

$ws.Cells.Item(38, 8).Value = 776.2727
$ws.Cells.Item(38, 9).Value = 776.2727
$ws.Cells.Item(38, 11).Value = 2328.8181
$ws.Cells.Item(38, 13).Value = -1956.8181

$ws.Cells.Item(74, 8).Value = 3917.5
$ws.Cells.Item(74, 9).Value = 3701
$ws.Cells.Item(74, 11).Value = 3701
$ws.Cells.Item(74, 13).Value = -2765

$ws.Cells.Item(77, 8).Value = 3917.5
$ws.Cells.Item(77, 9).Value = 3701
$ws.Cells.Item(77, 11).Value = 18505
$ws.Cells.Item(77, 13).Value = -13825

$ws.Cells.Item(80, 8).Value = 1508.6428
$ws.Cells.Item(80, 10).Value = 966.5
$ws.Cells.Item(80, 12).Value = 2899.5
$ws.Cells.Item(80, 14).Value = -4895.5

$ws.Cells.Item(83, 8).Value = 1508.6428
$ws.Cells.Item(83, 10).Value = 966.5
$ws.Cells.Item(83, 12).Value = 8698.5
$ws.Cells.Item(83, 14).Value = -18682.5

$ws.Cells.Item(88, 8).Value = 4133.5
$ws.Cells.Item(88, 9).Value = 3998.5
$ws.Cells.Item(88, 10).Value = 4160.5
$ws.Cells.Item(88, 11).Value = 3998.5
$ws.Cells.Item(88, 12).Value = 4160.5
$ws.Cells.Item(88, 13).Value = -3592.5
$ws.Cells.Item(88, 14).Value = -4972.5

$ws.Cells.Item(91, 8).Value = 4133.5
$ws.Cells.Item(91, 9).Value = 3998.5
$ws.Cells.Item(91, 10).Value = 4160.5
$ws.Cells.Item(91, 11).Value = 3998.5
$ws.Cells.Item(91, 12).Value = 4160.5
$ws.Cells.Item(91, 13).Value = -2594.5
$ws.Cells.Item(91, 14).Value = -6968.5

$ws.Cells.Item(100, 8).Value = 1275
$ws.Cells.Item(100, 10).Value = 1064.3334
$ws.Cells.Item(100, 12).Value = 1064.3334
$ws.Cells.Item(100, 14).Value = -2146.3334

$ws.Cells.Item(112, 8).Value = 1925.5385
$ws.Cells.Item(112, 10).Value = 2116.5454
$ws.Cells.Item(112, 12).Value = 6349.6362
$ws.Cells.Item(112, 14).Value = -8565.636200000001

$ws.Cells.Item(138, 8).Value = 4215.635
$ws.Cells.Item(138, 10).Value = 4998
$ws.Cells.Item(138, 12).Value = 14994
$ws.Cells.Item(138, 14).Value = -25274

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(36, 8).Value = 12000
$ws.Cells.Item(36, 9).Value = 8000
$ws.Cells.Item(36, 11).Value = 8000
$ws.Cells.Item(36, 13).Value = -7654

$ws.Cells.Item(96, 8).Value = 12499.5
$ws.Cells.Item(96, 10).Value = 12499.5
$ws.Cells.Item(96, 12).Value = 12499.5
$ws.Cells.Item(96, 14).Value = -17991.5

$ws.Cells.Item(101, 8).Value = 32500
$ws.Cells.Item(101, 9).Value = 30000
$ws.Cells.Item(101, 11).Value = 30000
$ws.Cells.Item(101, 13).Value = -26755

$ws.Cells.Item(110, 8).Value = 6317.1816
$ws.Cells.Item(110, 9).Value = 5944.3335
$ws.Cells.Item(110, 10).Value = 7995
$ws.Cells.Item(110, 11).Value = 5944.3335
$ws.Cells.Item(110, 12).Value = 7995
$ws.Cells.Item(110, 13).Value = -3899.3335
$ws.Cells.Item(110, 14).Value = -12085

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 1812.2
$ws.Cells.Item(20, 9).Value = 1612.5555
$ws.Cells.Item(20, 11).Value = 1612.5555
$ws.Cells.Item(20, 13).Value = -1365.5555

$ws.Cells.Item(23, 8).Value = 980
$ws.Cells.Item(23, 9).Value = 980
$ws.Cells.Item(23, 10).Value = 0
$ws.Cells.Item(23, 11).Value = 980
$ws.Cells.Item(23, 12).Value = 0
$ws.Cells.Item(23, 13).Value = -697
$ws.Cells.Item(23, 14).ClearContents()

$ws.Cells.Item(31, 8).Value = 6745.75
$ws.Cells.Item(31, 9).Value = 1327.6666
$ws.Cells.Item(31, 11).Value = 1327.6666
$ws.Cells.Item(31, 13).Value = -1075.6666

$ws.Cells.Item(86, 8).Value = 2333
$ws.Cells.Item(86, 9).Value = 2998
$ws.Cells.Item(86, 11).Value = 2998
$ws.Cells.Item(86, 13).Value = -1875

$ws.Cells.Item(89, 8).Value = 2333
$ws.Cells.Item(89, 9).Value = 2998
$ws.Cells.Item(89, 11).Value = 14990
$ws.Cells.Item(89, 13).Value = -9374

$ws.Cells.Item(107, 8).Value = 2900
$ws.Cells.Item(107, 9).Value = 2900
$ws.Cells.Item(107, 11).Value = 2900
$ws.Cells.Item(107, 13).Value = -980

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 104.045456
$ws.Cells.Item(7, 9).Value = 58.52941
$ws.Cells.Item(7, 10).Value = 258.8
$ws.Cells.Item(7, 11).Value = 58.52941
$ws.Cells.Item(7, 12).Value = 258.8
$ws.Cells.Item(7, 13).Value = 54.47059
$ws.Cells.Item(7, 14).Value = -484.8

$ws.Cells.Item(31, 8).Value = 3910
$ws.Cells.Item(31, 9).Value = 3441.1765
$ws.Cells.Item(31, 11).Value = 3441.1765
$ws.Cells.Item(31, 13).Value = -3146.1765

$ws.Cells.Item(34, 8).Value = 3910
$ws.Cells.Item(34, 9).Value = 3441.1765
$ws.Cells.Item(34, 11).Value = 3441.1765
$ws.Cells.Item(34, 13).Value = -3239.1765

$ws.Cells.Item(35, 8).Value = 8028
$ws.Cells.Item(35, 10).Value = 8028
$ws.Cells.Item(35, 12).Value = 8028
$ws.Cells.Item(35, 14).Value = -8616

$ws.Cells.Item(107, 8).Value = 1532.3334
$ws.Cells.Item(107, 9).Value = 823.5
$ws.Cells.Item(107, 11).Value = 823.5
$ws.Cells.Item(107, 13).Value = 1096.5

$ws.Cells.Item(122, 8).Value = 3787
$ws.Cells.Item(122, 9).Value = 4382.375
$ws.Cells.Item(122, 10).Value = 2199.3333
$ws.Cells.Item(122, 11).Value = 13147.125
$ws.Cells.Item(122, 12).Value = 6597.999899999999
$ws.Cells.Item(122, 13).Value = -10697.125
$ws.Cells.Item(122, 14).Value = -11497.9999

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(17, 8).Value = 291.8
$ws.Cells.Item(17, 9).Value = 190
$ws.Cells.Item(17, 10).Value = 444.5
$ws.Cells.Item(17, 11).Value = 570
$ws.Cells.Item(17, 12).Value = 1333.5
$ws.Cells.Item(17, 13).Value = -401
$ws.Cells.Item(17, 14).Value = -1671.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 2610.4285
$ws.Cells.Item(102, 9).Value = 2562.1667
$ws.Cells.Item(102, 11).Value = 2562.1667
$ws.Cells.Item(102, 13).Value = -940.1667000000002

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(35, 8).Value = 506.83334
$ws.Cells.Item(35, 9).Value = 506.83334
$ws.Cells.Item(35, 11).Value = 506.83334
$ws.Cells.Item(35, 13).Value = -170.83334

$ws.Cells.Item(93, 8).Value = 1449.9166
$ws.Cells.Item(93, 9).Value = 1355.4445
$ws.Cells.Item(93, 11).Value = 1355.4445
$ws.Cells.Item(93, 13).Value = -107.4445000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(47, 8).Value = 0
$ws.Cells.Item(47, 10).Value = 0
$ws.Cells.Item(47, 12).Value = 0
$ws.Cells.Item(47, 14).ClearContents()

$ws.Cells.Item(107, 8).Value = 1343.2222
$ws.Cells.Item(107, 9).Value = 1290.6666
$ws.Cells.Item(107, 10).Value = 1448.3334
$ws.Cells.Item(107, 11).Value = 3871.9998
$ws.Cells.Item(107, 12).Value = 4345.0002
$ws.Cells.Item(107, 13).Value = -1951.9998
$ws.Cells.Item(107, 14).Value = -8185.0002

$ws.Cells.Item(113, 8).Value = 693.55
$ws.Cells.Item(113, 9).Value = 631.3333
$ws.Cells.Item(113, 11).Value = 1893.9999
$ws.Cells.Item(113, 13).Value = 276.0001

$ws.Cells.Item(126, 8).Value = 4708.25
$ws.Cells.Item(126, 9).Value = 4094
$ws.Cells.Item(126, 10).Value = 5936.75
$ws.Cells.Item(126, 11).Value = 12282
$ws.Cells.Item(126, 12).Value = 17810.25
$ws.Cells.Item(126, 13).Value = -9812
$ws.Cells.Item(126, 14).Value = -22750.25

$ws.Cells.Item(136, 8).Value = 2220.8
$ws.Cells.Item(136, 9).Value = 1776
$ws.Cells.Item(136, 10).Value = 4000
$ws.Cells.Item(136, 11).Value = 5328
$ws.Cells.Item(136, 12).Value = 12000
$ws.Cells.Item(136, 13).Value = -2778
$ws.Cells.Item(136, 14).Value = -17100
